$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order corresponding to each row of data below: A(row idx), B, C, D, E, F, I, J, K, L, M, N
$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
    @(2,1.02,1.035511682083039,1.03801135020667,1.05317367678067,1.058945210776702,1.03864372796957,1.04062505500605,1.040800599443739,1.055920375723657,1.061676059821212,1.0174845810174),
    @(3,1.02,1.036482519387162,1.038731259459284,1.054458096582255,1.060307389221325,1.038911985225507,1.041239363595641,1.041330672385524,1.057016579847988,1.062850983470368,1.017689361002727),
    @(4,1.02,1.037110368941933,1.039196727220969,1.055289787392761,1.06118938298625,1.039083997356095,1.04163591361019,1.041672634513445,1.057725929407401,1.063611272047857,1.017821511258822),
    @(5,1.02,1.037374233940513,1.039392322754204,1.055639572270549,1.061560312424626,1.039155935671901,1.041802396034907,1.041816148461871,1.058024149354668,1.063930907301618,1.0178769819035),
    @(6,1.02,1.037418533174451,1.039425158980966,1.055698311077561,1.061622601334797,1.039167992403105,1.041830335841113,1.041840230598621,1.058074222395921,1.063984576099245,1.017886290663243),
    @(7,1.02,1.037113895044737,1.039199341121034,1.055294460674783,1.061194338813821,1.039084960075362,1.041638139049348,1.041674553124406,1.057729914197912,1.063615542988595,1.017822252795649),
    @(8,1.02,1.035839853051306,1.038254721477007,1.053607632437742,1.05940544788962,1.038734711398074,1.040832859864835,1.040979953317112,1.056290837388673,1.062073124665259,1.017553861058028),
    @(9,1.02,1.033592162828015,1.036587433868439,1.050639621684401,1.056257492898296,1.038105523401001,1.039406591704153,1.03974809655138,1.053755171218597,1.05935537397483,1.017078196417213),
    @(10,1.02,1.032091897317371,1.035474090230713,1.048663780542634,1.054161617861556,1.037677999562609,1.038450863863562,1.038921568871639,1.05206473482828,1.057543557850101,1.016759258462094),
    @(11,1.02,1.031441833579004,1.034991574057097,1.047808863438217,1.053254706497086,1.037490964898982,1.038035864494887,1.038562420045498,1.051332734623532,1.056759000906156,1.016620721614795),
    @(12,1.02,1.031200304560227,1.034812281644463,1.047491402099101,1.052917929223464,1.037421204328413,1.037881540380285,1.038428827372134,1.051060830453608,1.056467575008337,1.016569197566155),
    @(13,1.02,1.031252116367025,1.034850743385495,1.04755949446824,1.052990165061738,1.037436181220868,1.037914651342602,1.038457491997684,1.051119155159339,1.056530087183212,1.016580252604258),
    @(14,1.02,1.031421870072265,1.034976755002603,1.047782620103349,1.053226866561994,1.037485204334832,1.038023111583679,1.038551381084947,1.051310259089347,1.056734911696301,1.016616463952871),
    @(15,1.02,1.031526452119678,1.035054386392241,1.047920107412136,1.053372718014444,1.037515370987356,1.038089914340926,1.038609204232082,1.0514280035746,1.056861109909002,1.016638766307594),
    @(16,1.02,1.032135030620477,1.035506104187516,1.048720531744992,1.054221819344882,1.037690372110335,1.038478381490426,1.038945377917466,1.052113314450958,1.057595625518544,1.016768443539713),
    @(17,1.02,1.032516658009197,1.035789339666351,1.049222785075882,1.054754601349137,1.037799633248305,1.038721745217267,1.039155914303624,1.052543182262806,1.058056358530856,1.016849670284886),
    @(18,1.02,1.032739212495185,1.035954504518771,1.049515802395726,1.055065423723171,1.037863178747045,1.038863582968922,1.039278595342957,1.052793914127459,1.058325093690798,1.016897006537836),
    @(19,1.02,1.032815090660655,1.036010814397029,1.049615724228008,1.055171416316626,1.037884814796191,1.038911926951699,1.039320405836226,1.052879406812465,1.058416725060177,1.016913139864186),
    @(20,1.02,1.032475717381931,1.035758955481412,1.049168891708415,1.054697432731488,1.037787929663978,1.038695646198195,1.039133338285408,1.052497061815501,1.058006926524799,1.016840959761168),
    @(21,1.02,1.031371883673866,1.034939649497812,1.047716912619471,1.053157161406701,1.037470776188026,1.037991177574571,1.038523738300727,1.051253983947608,1.056674596180022,1.016605802415726),
    @(22,1.02,1.030677475193785,1.034424146744938,1.0468045322206,1.05218925103583,1.037269705419663,1.037547238358539,1.038139366218159,1.050472371701362,1.055836868410807,1.016457571880092),
    @(23,1.02,1.031045630862941,1.034697459712205,1.047288152310816,1.052702310282385,1.037376454550034,1.037782674900377,1.038343232669425,1.05088672364776,1.056280968170949,1.016536187501283),
    @(24,1.02,1.032494216819628,1.035772684915216,1.049193243610864,1.05472326458405,1.037793218580491,1.038707439568001,1.03914353979087,1.052517901668364,1.058029262728045,1.016844895801421),
    @(25,1.02,1.034173561845715,1.037018790269218,1.051406415312721,1.057070818508356,1.038269605577242,1.039776176827122,1.040067494161296,1.054410692483306,1.060057966616152,1.017201489598143)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $c = $colIndexes[$i]
        $val = $row[$i + 1]
        $ws.Cells.Item($r, $c).Value = $val
    }
}

Write-Host "Updated vm_pu values for rows 2-25"
